$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4335.3335
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4335.3335
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4335.3335
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4831.3335

$ws.Range("H67").Value = 4335.3335
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4335.3335
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4335.3335
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6051.3335

$ws.Range("H101").Value = 2317
$ws.Range("I101").Value = 2307.2
$ws.Range("J101").Value = 2333.3333
$ws.Range("K101").Value = 6921.599999999999
$ws.Range("L101").Value = 6999.999899999999
$ws.Range("M101").Value = -5299.599999999999
$ws.Range("N101").Value = -10243.9999

$ws.Range("H129").Value = 2021.08
$ws.Range("I129").Value = 694.375
$ws.Range("J129").Value = 2645.4119
$ws.Range("K129").Value = 2083.125
$ws.Range("L129").Value = 7936.2357
$ws.Range("M129").Value = 2916.875
$ws.Range("N129").Value = -17936.2357

$ws.Range("H137").Value = 1840.5416
$ws.Range("I137").Value = 2378.111
$ws.Range("J137").Value = 1518
$ws.Range("K137").Value = 7134.333
$ws.Range("L137").Value = 4554
$ws.Range("M137").Value = -4584.333
$ws.Range("N137").Value = -9654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2002.2858
$ws.Range("I61").Value = 2121.3333
$ws.Range("J61").Value = 1913
$ws.Range("K61").Value = 2121.3333
$ws.Range("L61").Value = 1913
$ws.Range("M61").Value = -1909.3333
$ws.Range("N61").Value = -2337

$ws.Range("H122").Value = 1645.8235
$ws.Range("I122").Value = 2033.3334
$ws.Range("J122").Value = 1434.4546
$ws.Range("K122").Value = 6100.0002
$ws.Range("L122").Value = 4303.3638
$ws.Range("M122").Value = -3650.0002
$ws.Range("N122").Value = -9203.363799999999

$ws.Range("H136").Value = 2002.2858
$ws.Range("I136").Value = 2121.3333
$ws.Range("J136").Value = 1913
$ws.Range("K136").Value = 6363.999899999999
$ws.Range("L136").Value = 5739
$ws.Range("M136").Value = -3813.999899999999
$ws.Range("N136").Value = -10839

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 99999
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99999
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99999
$ws.Range("N140").Value = -110359

$ws.Range("H141").Value = 69285.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 69285.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 69285.8
$ws.Range("N141").Value = -79645.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 9500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 9500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 9500
$ws.Range("N76").Value = -10130

$ws.Range("H79").Value = 9500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 9500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 9500
$ws.Range("N79").Value = -11684

$ws.Range("H134").Value = 2229.1804
$ws.Range("I134").Value = 1455.0731
$ws.Range("J134").Value = 3816.1
$ws.Range("K134").Value = 4365.219300000001
$ws.Range("L134").Value = 11448.3
$ws.Range("M134").Value = -1830.219300000001
$ws.Range("N134").Value = -16518.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2267.8171
$ws.Range("I31").Value = 1181.898
$ws.Range("J31").Value = 3477.1365
$ws.Range("K31").Value = 1181.898
$ws.Range("L31").Value = 3477.1365
$ws.Range("M31").Value = -886.8979999999999
$ws.Range("N31").Value = -4067.1365

$ws.Range("H34").Value = 2267.8171
$ws.Range("I34").Value = 1181.898
$ws.Range("J34").Value = 3477.1365
$ws.Range("K34").Value = 1181.898
$ws.Range("L34").Value = 3477.1365
$ws.Range("M34").Value = -979.8979999999999
$ws.Range("N34").Value = -3881.1365

$ws.Range("H62").Value = 10996.917
$ws.Range("I62").Value = 2938.75
$ws.Range("J62").Value = 15026
$ws.Range("K62").Value = 2938.75
$ws.Range("L62").Value = 15026
$ws.Range("M62").Value = -2314.75
$ws.Range("N62").Value = -16274

$ws.Range("H65").Value = 10996.917
$ws.Range("I65").Value = 2938.75
$ws.Range("J65").Value = 15026
$ws.Range("K65").Value = 14693.75
$ws.Range("L65").Value = 75130
$ws.Range("M65").Value = -11573.75
$ws.Range("N65").Value = -81370

$ws.Range("H105").Value = 2005
$ws.Range("I105").Value = 2005.5555
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2005.5555
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -258.5554999999999
$ws.Range("N105").Value = -5494

$ws.Range("H132").Value = 1576.5646
$ws.Range("I132").Value = 1049.3256
$ws.Range("J132").Value = 2769.7896
$ws.Range("K132").Value = 3147.976799999999
$ws.Range("L132").Value = 8309.3688
$ws.Range("M132").Value = -617.9767999999995
$ws.Range("N132").Value = -13369.3688

$ws.Range("H134").Value = 1815.4
$ws.Range("I134").Value = 1311.75
$ws.Range("J134").Value = 2391
$ws.Range("K134").Value = 3935.25
$ws.Range("L134").Value = 7173
$ws.Range("M134").Value = -1400.25
$ws.Range("N134").Value = -12243

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 576.94446
$ws.Range("I5").Value = 351.32257
$ws.Range("J5").Value = 1975.8
$ws.Range("K5").Value = 1053.96771
$ws.Range("L5").Value = 5927.4
$ws.Range("M5").Value = -941.9677099999999
$ws.Range("N5").Value = -6151.4

$ws.Range("H33").Value = 211
$ws.Range("I33").Value = 96.333336
$ws.Range("J33").Value = 260.14285
$ws.Range("K33").Value = 578.000016
$ws.Range("L33").Value = 1560.8571
$ws.Range("M33").Value = -295.000016
$ws.Range("N33").Value = -2126.8571

$ws.Range("H135").Value = 576.94446
$ws.Range("I135").Value = 351.32257
$ws.Range("J135").Value = 1975.8
$ws.Range("K135").Value = 3161.90313
$ws.Range("L135").Value = 17782.2
$ws.Range("M135").Value = -626.9031299999997
$ws.Range("N135").Value = -22852.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2073.2307
$ws.Range("I113").Value = 1848.1428
$ws.Range("J113").Value = 3018.6
$ws.Range("K113").Value = 1848.1428
$ws.Range("L113").Value = 3018.6
$ws.Range("M113").Value = 321.8571999999999
$ws.Range("N113").Value = -7358.6

$ws.Range("H126").Value = 4466669.5
$ws.Range("I126").Value = 7814493.5
$ws.Range("J126").Value = 2904.5
$ws.Range("K126").Value = 23443480.5
$ws.Range("L126").Value = 8713.5
$ws.Range("M126").Value = -23441010.5
$ws.Range("N126").Value = -13653.5

$ws.Range("H132").Value = 2545.12
$ws.Range("I132").Value = 1534.7037
$ws.Range("J132").Value = 3731.261
$ws.Range("K132").Value = 4604.1111
$ws.Range("L132").Value = 11193.783
$ws.Range("M132").Value = -2074.1111
$ws.Range("N132").Value = -16253.783

$ws.Range("H133").Value = 49813.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49813.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49813.332
$ws.Range("N133").Value = -59933.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4335
$ws.Range("I7").Value = 4500
$ws.Range("J7").Value = 4252.5
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 4252.5
$ws.Range("M7").Value = -4388
$ws.Range("N7").Value = -4476.5

$ws.Range("H64").Value = 27749.75
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 27749.75
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 27749.75
$ws.Range("N64").Value = -28199.75

$ws.Range("H67").Value = 27749.75
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 27749.75
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 27749.75
$ws.Range("N67").Value = -29309.75

$ws.Range("H122").Value = 8638.333000000001
$ws.Range("I122").Value = 17041.143
$ws.Range("J122").Value = 3291.0908
$ws.Range("K122").Value = 51123.429
$ws.Range("L122").Value = 9873.2724
$ws.Range("M122").Value = -48673.429
$ws.Range("N122").Value = -14773.2724

$ws.Range("H126").Value = 4335
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 4252.5
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 12757.5
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -17697.5

$ws.Range("H133").Value = 39800
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 39800
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 39800
$ws.Range("N133").Value = -44860

$ws.Range("H134").Value = 48500
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 48500
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 48500
$ws.Range("N134").Value = -58640

$ws.Range("H135").Value = 47000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 47000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 47000
$ws.Range("N135").Value = -57140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 30999.666
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 30999.666
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 30999.666
$ws.Range("N68").Value = -32621.666

$ws.Range("H69").Value = 9531.111000000001
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9531.111000000001
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9531.111000000001
$ws.Range("N69").Value = -11029.111

$ws.Range("H71").Value = 30999.666
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 30999.666
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 92998.99800000001
$ws.Range("N71").Value = -101110.998

$ws.Range("H72").Value = 9531.111000000001
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9531.111000000001
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 28593.333
$ws.Range("N72").Value = -36081.333

$ws.Range("H80").Value = 9499.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9499.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9499.5
$ws.Range("N80").Value = -11495.5

$ws.Range("H83").Value = 9499.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9499.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 28498.5
$ws.Range("N83").Value = -38482.5
